$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values (avoid Excel auto-converting to numbers)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '54.028.61'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.242.44'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '495.14'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '127.12'
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").Value = '2.283.19'
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("D10").Value = '0.0943'
$ws.Range("E10").Value = '  +3.31%  '
$ws.Range("E11").Value = '  +2.32%  '
$ws.Range("E12").Value = '  +3.14%  '
$ws.Range("D13").Value = '4.61'
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").Value = '2.670.27'
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").Value = '21.70'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '54.091.94'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '2.286.92'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = '10.01'
$ws.Range("E19").Value = '  +4.51%  '
$ws.Range("D20").Value = '4.08'
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("E21").Value = '  +6.20%  '
$ws.Range("D22").Value = '301.28'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("D25").Value = '62.29'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +1.99%  '
$ws.Range("D28").Value = '2.392.13'
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("E29").Value = '  +4.87%  '
$ws.Range("D30").Value = '7.08'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = '168.37'
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = '0.0₃0687'
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").Value = '1.59'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = '5.85'
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("D36").Value = '0.990'
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("E39").Value = '  +2.41%  '
$ws.Range("D40").Value = '0.867'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = '3.69'
$ws.Range("E41").Value = '  +3.15%  '
$ws.Range("D42").Value = '35.37'
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = '0.374'
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.40'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").Value = '128.48'
$ws.Range("D47").Value = '4.91'
$ws.Range("E47").Value = '  +5.10%  '
$ws.Range("D48").Value = '0.0888'
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").Value = '0.542'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").Value = '237.38'
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '0.0483'
$ws.Range("E51").Value = '  +2.59%  '
